$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 171, pushing the existing
# rows 171:189 down to 173:191 (formatting/styles of row 171 carry down
# with the insert, same as Excel's native "Insert Copied/Blank Rows").
$ws.Rows("171:172").Insert()

# New row 171: Vega Central Mapocho de Santiago, Berenjena, Primera,
# week of 2021-11-04 (serial 44504)
$ws.Range("A171").Value = 9
$ws.Range("B171").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C171").Value = "Metropolitana"
$ws.Range("D171").Value = 44504
$ws.Range("E171").Value = 13
$ws.Range("F171").Value = 100112001
$ws.Range("G171").Value = "Berenjena"
$ws.Range("H171").Value = "Sin especificar"
$ws.Range("I171").Value = "Primera"
$ws.Range("J171").Value = 70
$ws.Range("K171").Value = 9000
$ws.Range("L171").Value = 10000
$ws.Range("M171").Value = 9500
$ws.Range("N171").Value = "`$/caja 60 unidades"
$ws.Range("O171").Value = "Región de Arica y Parinacota"
$ws.Range("P171").Value = 158
$ws.Range("Q171").Value = 60
$ws.Range("R171").Value = "Hortaliza"

# New row 172: same market/date, Segunda quality
$ws.Range("A172").Value = 9
$ws.Range("B172").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C172").Value = "Metropolitana"
$ws.Range("D172").Value = 44504
$ws.Range("E172").Value = 13
$ws.Range("F172").Value = 100112001
$ws.Range("G172").Value = "Berenjena"
$ws.Range("H172").Value = "Sin especificar"
$ws.Range("I172").Value = "Segunda"
$ws.Range("J172").Value = 34
$ws.Range("K172").Value = 7000
$ws.Range("L172").Value = 7000
$ws.Range("M172").Value = 7000
$ws.Range("N172").Value = "`$/caja 100 unidades"
$ws.Range("O172").Value = "Región de Arica y Parinacota"
$ws.Range("P172").Value = 70
$ws.Range("Q172").Value = 100
$ws.Range("R172").Value = "Hortaliza"
